# Applies the changes described by the commit:
#   - rename "volume" -> "volume (ul)" and "concentration" -> "concentration (ng-ul)"
#   - make the "volume" sheet the active tab (was "concentration")
#   - update per-sheet selection / scroll position
#   - convert every 5E-005 placeholder value on the "volume" sheet to 50
#   - touch two extra (empty) trailing rows on the "volume" sheet

$wb = $excel.ActiveWorkbook

$wsContent = $wb.Worksheets.Item(1)
$wsVolume = $wb.Worksheets.Item(2)
$wsConcentration = $wb.Worksheets.Item(3)

# --- Rename sheets -----------------------------------------------------
$wsVolume.Name = "volume (ul)"
$wsConcentration.Name = "concentration (ng-ul)"

# --- Fix up the placeholder volume values (5E-005 -> 50) ---------------
$used = $wsVolume.UsedRange
$rowCount = $used.Rows.Count
$colCount = $used.Columns.Count
for ($r = 1; $r -le $rowCount; $r++) {
  for ($c = 1; $c -le $colCount; $c++) {
    $cell = $wsVolume.Cells.Item($r, $c)
    if ($cell.Value2 -eq 0.00005) {
      $cell.Value = 50
    }
  }
}

# --- Touch two extra trailing (empty) rows on the volume sheet ---------
$defaultRowHeight = $wsVolume.Rows.Item(18).RowHeight
$wsVolume.Rows.Item(19).RowHeight = $defaultRowHeight
$wsVolume.Rows.Item(20).RowHeight = $defaultRowHeight

# --- Update selection / scroll position on each sheet -------------------
$null = $wsContent.Activate()
$null = $excel.Goto($wsContent.Range("E1"), $true)
$null = $wsContent.Range("Y16").Select()

$null = $wsConcentration.Activate()
$null = $excel.Goto($wsConcentration.Range("H1"), $true)
$null = $wsConcentration.Range("Y16").Select()

# volume becomes the active tab (activeTab=1), matches the new selection Y17
$null = $wsVolume.Activate()
$null = $excel.Goto($wsVolume.Range("J1"), $true)
$null = $wsVolume.Range("Y17").Select()
